# Week 22 profile updates
#
# 1. Flip the "completed" flag for row 10 (Manas / manas_rishav) from 0 to 1.
# 2. Append a brand-new profile row (row 97): Ashutosh Menghrajani / ashu_menghrajani,
#    with the flag column starting at 0 - mirrors the pattern used by every other
#    row in the sheet (name, handle, 0/1 flag).
# 3. Reflect where the author had scrolled to / was last looking when they saved -
#    selection on F96, window scrolled down near row 88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. C10: 0 -> 1 ------------------------------------------------------
$ws.Range("C10").Value = 1

# --- 2. New row 97 --------------------------------------------------------
$ws.Range("A97").Value = "Ashutosh Menghrajani"
$ws.Range("B97").Value = "ashu_menghrajani"

# Give the new flag cell the same look used by every other column-C status
# flag in this sheet (copy format from an existing flag cell, then set the
# value for the new profile).
$ws.Range("C10").Copy()
$ws.Range("C97").PasteSpecial(-4122)
$ws.Range("C97").Value = 0

# --- 3. Scroll position / selection --------------------------------------
$excel.ActiveWindow.ScrollRow = 88
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F96").Select()

Write-Host "Week 22 profile updates applied"
